$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-12-25 Monday", $true, $true, $false, $false, $false, $true, 1, $false, "2023-12-26 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("64-13=51", $true, $true, $false, $false, $false, $true, 1, $false, "86-41=45", 2) | Out-Null
$d.Content.Find.Execute("94-56=38", $true, $true, $false, $false, $false, $true, 1, $false, "31+15=46", 2) | Out-Null
$d.Content.Find.Execute("26-13=13", $true, $true, $false, $false, $false, $true, 1, $false, "43-3=40", 2) | Out-Null
$d.Content.Find.Execute("40+29=69", $true, $true, $false, $false, $false, $true, 1, $false, "59-52=7", 2) | Out-Null
$d.Content.Find.Execute("32-23=9", $true, $true, $false, $false, $false, $true, 1, $false, "51-5=46", 2) | Out-Null
$d.Content.Find.Execute("75+2=77", $true, $true, $false, $false, $false, $true, 1, $false, "46-22=24", 2) | Out-Null
$d.Content.Find.Execute("17+3=20", $true, $true, $false, $false, $false, $true, 1, $false, "8+80=88", 2) | Out-Null
$d.Content.Find.Execute("79-48=31", $true, $true, $false, $false, $false, $true, 1, $false, "22-15=7", 2) | Out-Null
$d.Content.Find.Execute("4+54=58", $true, $true, $false, $false, $false, $true, 1, $false, "52-22=30", 2) | Out-Null
$d.Content.Find.Execute("90+2=92", $true, $true, $false, $false, $false, $true, 1, $false, "46-21=25", 2) | Out-Null
$d.Content.Find.Execute("24+73=97", $true, $true, $false, $false, $false, $true, 1, $false, "12+17=29", 2) | Out-Null
$d.Content.Find.Execute("95-27=68", $true, $true, $false, $false, $false, $true, 1, $false, "62-8=54", 2) | Out-Null
$d.Content.Find.Execute("3+52=55", $true, $true, $false, $false, $false, $true, 1, $false, "22+9=31", 2) | Out-Null
$d.Content.Find.Execute("47+37=84", $true, $true, $false, $false, $false, $true, 1, $false, "28+41=69", 2) | Out-Null
$d.Content.Find.Execute("77+21=98", $true, $true, $false, $false, $false, $true, 1, $false, "81-49=32", 2) | Out-Null
$d.Content.Find.Execute("53+30=83", $true, $true, $false, $false, $false, $true, 1, $false, "5+92=97", 2) | Out-Null
$d.Content.Find.Execute("56+9=65", $true, $true, $false, $false, $false, $true, 1, $false, "46-2=44", 2) | Out-Null
$d.Content.Find.Execute("54-15=39", $true, $true, $false, $false, $false, $true, 1, $false, "67-0=67", 2) | Out-Null
$d.Content.Find.Execute("54-35=19", $true, $true, $false, $false, $false, $true, 1, $false, "95-61=34", 2) | Out-Null
$d.Content.Find.Execute("12+25=37", $true, $true, $false, $false, $false, $true, 1, $false, "62-56=6", 2) | Out-Null
$d.Content.Find.Execute("73-1=72", $true, $true, $false, $false, $false, $true, 1, $false, "86+6=92", 2) | Out-Null
$d.Content.Find.Execute("74-68=6", $true, $true, $false, $false, $false, $true, 1, $false, "18+61=79", 2) | Out-Null
$d.Content.Find.Execute("32+58=90", $true, $true, $false, $false, $false, $true, 1, $false, "92-33=59", 2) | Out-Null
$d.Content.Find.Execute("95-68=27", $true, $true, $false, $false, $false, $true, 1, $false, "32-32=0", 2) | Out-Null
$d.Content.Find.Execute("42-9=33", $true, $true, $false, $false, $false, $true, 1, $false, "90+7=97", 2) | Out-Null
$d.Content.Find.Execute("57+6=63", $true, $true, $false, $false, $false, $true, 1, $false, "50+28=78", 2) | Out-Null
$d.Content.Find.Execute("67-47=20", $true, $true, $false, $false, $false, $true, 1, $false, "95-50=45", 2) | Out-Null
$d.Content.Find.Execute("35-19=16", $true, $true, $false, $false, $false, $true, 1, $false, "22+33=55", 2) | Out-Null
$d.Content.Find.Execute("36+60=96", $true, $true, $false, $false, $false, $true, 1, $false, "77-7=70", 2) | Out-Null
$d.Content.Find.Execute("9+31=40", $true, $true, $false, $false, $false, $true, 1, $false, "33+56=89", 2) | Out-Null
$d.Content.Find.Execute("15+67=82", $true, $true, $false, $false, $false, $true, 1, $false, "23-12=11", 2) | Out-Null
$d.Content.Find.Execute("22+37=59", $true, $true, $false, $false, $false, $true, 1, $false, "84-58=26", 2) | Out-Null
$d.Content.Find.Execute("79+19=98", $true, $true, $false, $false, $false, $true, 1, $false, "61+18=79", 2) | Out-Null
$d.Content.Find.Execute("36-12=24", $true, $true, $false, $false, $false, $true, 1, $false, "72-44=28", 2) | Out-Null
$d.Content.Find.Execute("78-44=34", $true, $true, $false, $false, $false, $true, 1, $false, "47+29=76", 2) | Out-Null
$d.Content.Find.Execute("39+6=45", $true, $true, $false, $false, $false, $true, 1, $false, "93-86=7", 2) | Out-Null
$d.Content.Find.Execute("11+48=59", $true, $true, $false, $false, $false, $true, 1, $false, "86-48=38", 2) | Out-Null
$d.Content.Find.Execute("86-59=27", $true, $true, $false, $false, $false, $true, 1, $false, "1+63=64", 2) | Out-Null
$d.Content.Find.Execute("20-3=17", $true, $true, $false, $false, $false, $true, 1, $false, "44-32=12", 2) | Out-Null
$d.Content.Find.Execute("64-26=38", $true, $true, $false, $false, $false, $true, 1, $false, "39+52=91", 2) | Out-Null
$d.Content.Find.Execute("84-42=42", $true, $true, $false, $false, $false, $true, 1, $false, "71-38=33", 2) | Out-Null
$d.Content.Find.Execute("65-39=26", $true, $true, $false, $false, $false, $true, 1, $false, "72-55=17", 2) | Out-Null
$d.Content.Find.Execute("7+28=35", $true, $true, $false, $false, $false, $true, 1, $false, "69-49=20", 2) | Out-Null
$d.Content.Find.Execute("62-26=36", $true, $true, $false, $false, $false, $true, 1, $false, "64-58=6", 2) | Out-Null
$d.Content.Find.Execute("16+27=43", $true, $true, $false, $false, $false, $true, 1, $false, "30-10=20", 2) | Out-Null
$d.Content.Find.Execute("72+8=80", $true, $true, $false, $false, $false, $true, 1, $false, "39-27=12", 2) | Out-Null
$d.Content.Find.Execute("12+29=41", $true, $true, $false, $false, $false, $true, 1, $false, "18+26=44", 2) | Out-Null
$d.Content.Find.Execute("68-16=52", $true, $true, $false, $false, $false, $true, 1, $false, "58-16=42", 2) | Out-Null
$d.Content.Find.Execute("87-0=87", $true, $true, $false, $false, $false, $true, 1, $false, "51-29=22", 2) | Out-Null
$d.Content.Find.Execute("78-14=64", $true, $true, $false, $false, $false, $true, 1, $false, "27-9=18", 2) | Out-Null
$d.Content.Find.Execute("31-14=17", $true, $true, $false, $false, $false, $true, 1, $false, "42+32=74", 2) | Out-Null
$d.Content.Find.Execute("90-15=75", $true, $true, $false, $false, $false, $true, 1, $false, "83-14=69", 2) | Out-Null
$d.Content.Find.Execute("30-22=8", $true, $true, $false, $false, $false, $true, 1, $false, "45+41=86", 2) | Out-Null
$d.Content.Find.Execute("37-11=26", $true, $true, $false, $false, $false, $true, 1, $false, "71-44=27", 2) | Out-Null
$d.Content.Find.Execute("74-57=17", $true, $true, $false, $false, $false, $true, 1, $false, "93-33=60", 2) | Out-Null
$d.Content.Find.Execute("93-92=1", $true, $true, $false, $false, $false, $true, 1, $false, "4+67=71", 2) | Out-Null
$d.Content.Find.Execute("69-29=40", $true, $true, $false, $false, $false, $true, 1, $false, "49+29=78", 2) | Out-Null
$d.Content.Find.Execute("47-38=9", $true, $true, $false, $false, $false, $true, 1, $false, "35+25=60", 2) | Out-Null
$d.Content.Find.Execute("35+36=71", $true, $true, $false, $false, $false, $true, 1, $false, "19+46=65", 2) | Out-Null
$d.Content.Find.Execute("10+64=74", $true, $true, $false, $false, $false, $true, 1, $false, "47+8=55", 2) | Out-Null
$d.Content.Find.Execute("53+6=59", $true, $true, $false, $false, $false, $true, 1, $false, "7+80=87", 2) | Out-Null
$d.Content.Find.Execute("37+38=75", $true, $true, $false, $false, $false, $true, 1, $false, "78-36=42", 2) | Out-Null
$d.Content.Find.Execute("48-1=47", $true, $true, $false, $false, $false, $true, 1, $false, "26-8=18", 2) | Out-Null
$d.Content.Find.Execute("59+24=83", $true, $true, $false, $false, $false, $true, 1, $false, "0+20=20", 2) | Out-Null
$d.Content.Find.Execute("37-24=13", $true, $true, $false, $false, $false, $true, 1, $false, "48-44=4", 2) | Out-Null
$d.Content.Find.Execute("51+20=71", $true, $true, $false, $false, $false, $true, 1, $false, "49-47=2", 2) | Out-Null
$d.Content.Find.Execute("52+30=82", $true, $true, $false, $false, $false, $true, 1, $false, "9-7=2", 2) | Out-Null
$d.Content.Find.Execute("48+50=98", $true, $true, $false, $false, $false, $true, 1, $false, "89-32=57", 2) | Out-Null
$d.Content.Find.Execute("23+25=48", $true, $true, $false, $false, $false, $true, 1, $false, "41-6=35", 2) | Out-Null
$d.Content.Find.Execute("94-6=88", $true, $true, $false, $false, $false, $true, 1, $false, "76-36=40", 2) | Out-Null
$d.Content.Find.Execute("11+57=68", $true, $true, $false, $false, $false, $true, 1, $false, "58-52=6", 2) | Out-Null
$d.Content.Find.Execute("35+33=68", $true, $true, $false, $false, $false, $true, 1, $false, "0+84=84", 2) | Out-Null
$d.Content.Find.Execute("34+35=69", $true, $true, $false, $false, $false, $true, 1, $false, "54-12=42", 2) | Out-Null
$d.Content.Find.Execute("1+79=80", $true, $true, $false, $false, $false, $true, 1, $false, "75+22=97", 2) | Out-Null
$d.Content.Find.Execute("33+53=86", $true, $true, $false, $false, $false, $true, 1, $false, "10+72=82", 2) | Out-Null
$d.Content.Find.Execute("52-38=14", $true, $true, $false, $false, $false, $true, 1, $false, "97-33=64", 2) | Out-Null
$d.Content.Find.Execute("42-29=13", $true, $true, $false, $false, $false, $true, 1, $false, "66-3=63", 2) | Out-Null
$d.Content.Find.Execute("34+54=88", $true, $true, $false, $false, $false, $true, 1, $false, "26+0=26", 2) | Out-Null
$d.Content.Find.Execute("72+15=87", $true, $true, $false, $false, $false, $true, 1, $false, "87-42=45", 2) | Out-Null
$d.Content.Find.Execute("99-29=70", $true, $true, $false, $false, $false, $true, 1, $false, "36+12=48", 2) | Out-Null
$d.Content.Find.Execute("81-22=59", $true, $true, $false, $false, $false, $true, 1, $false, "32+22=54", 2) | Out-Null
$d.Content.Find.Execute("28+28=56", $true, $true, $false, $false, $false, $true, 1, $false, "59-44=15", 2) | Out-Null
$d.Content.Find.Execute("88-32=56", $true, $true, $false, $false, $false, $true, 1, $false, "54+0=54", 2) | Out-Null
$d.Content.Find.Execute("36+16=52", $true, $true, $false, $false, $false, $true, 1, $false, "30-5=25", 2) | Out-Null
$d.Content.Find.Execute("21-6=15", $true, $true, $false, $false, $false, $true, 1, $false, "71-27=44", 2) | Out-Null
$d.Content.Find.Execute("92-88=4", $true, $true, $false, $false, $false, $true, 1, $false, "83-71=12", 2) | Out-Null
$d.Content.Find.Execute("79-70=9", $true, $true, $false, $false, $false, $true, 1, $false, "17+81=98", 2) | Out-Null
$d.Content.Find.Execute("80-39=41", $true, $true, $false, $false, $false, $true, 1, $false, "70+5=75", 2) | Out-Null
$d.Content.Find.Execute("76-21=55", $true, $true, $false, $false, $false, $true, 1, $false, "19+71=90", 2) | Out-Null
$d.Content.Find.Execute("55+28=83", $true, $true, $false, $false, $false, $true, 1, $false, "13+86=99", 2) | Out-Null
$d.Content.Find.Execute("66+9=75", $true, $true, $false, $false, $false, $true, 1, $false, "29-1=28", 2) | Out-Null
$d.Content.Find.Execute("98-34=64", $true, $true, $false, $false, $false, $true, 1, $false, "36+45=81", 2) | Out-Null
$d.Content.Find.Execute("40-20=20", $true, $true, $false, $false, $false, $true, 1, $false, "19+34=53", 2) | Out-Null
$d.Content.Find.Execute("34-15=19", $true, $true, $false, $false, $false, $true, 1, $false, "64+6=70", 2) | Out-Null
$d.Content.Find.Execute("54+9=63", $true, $true, $false, $false, $false, $true, 1, $false, "72+22=94", 2) | Out-Null
$d.Content.Find.Execute("13+49=62", $true, $true, $false, $false, $false, $true, 1, $false, "37-23=14", 2) | Out-Null
$d.Content.Find.Execute("18-0=18", $true, $true, $false, $false, $false, $true, 1, $false, "31+11=42", 2) | Out-Null
$d.Content.Find.Execute("44-7=37", $true, $true, $false, $false, $false, $true, 1, $false, "37+42=79", 2) | Out-Null
$d.Content.Find.Execute("59-5=54", $true, $true, $false, $false, $false, $true, 1, $false, "30+16=46", 2) | Out-Null
$d.Content.Find.Execute("79-3=76", $true, $true, $false, $false, $false, $true, 1, $false, "96-36=60", 2) | Out-Null
